# Update the "想去人数" (want-to-go count) column F for a handful of rows
# on both the "展览" and "全部类型" worksheets, incrementing each value by 1
# (refreshed data snapshot).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
$rows = @(4, 7, 17, 19, 22)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($r in $rows) {
        $cell = $ws.Cells.Item($r, 6)  # Column F
        $current = $cell.Value()
        $cell.Value = $current + 1
    }
}

$wb.Save()
